$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 8421.18
$ws.Range("B16").Value = 8046.99
$ws.Range("C16").Value = 17.2
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = $false
$ws.Range("F16").Value = 4.6500000000000004
$ws.Range("G16").Value = 42626.545624999999
$ws.Range("G16").NumberFormat = "m/d/yy h:mm"
$ws.Range("H16").Value = $true
